$d = $word.ActiveDocument

# 1) The first "Graph1" picture run gains a <w:rPr><w:noProof/></w:rPr>.
$d.InlineShapes(1).Range.NoProofing = 1

# 2) "Basic:" / " Exponential" -> "Basic:" / " " / "Polynomial" (Graph1 section)
$d.Content.Find.Execute("Basic: Exponential", $false, $false, $false, $false, $false, $true, 1, $false, "Basic: Polynomial", 2)

# 3) "exponentially" -> "polynomially" in the memory-growth explanation paragraph
$d.Content.Find.Execute("takes exponentially more memory", $false, $false, $false, $false, $false, $true, 1, $false, "takes polynomially more memory", 2)

# 4) "show an ... polynomial exponential increasing" -> "show a ... polynomial increasing" (Graph2 explanation)
$d.Content.Find.Execute("show an polynomial exponential increasing", $false, $false, $false, $false, $false, $true, 1, $false, "show a polynomial increasing", 2)

# 5) O(mn) -> O(m*n)
$d.Content.Find.Execute("Both algorithms take O(mn) time", $false, $false, $false, $false, $false, $true, 1, $false, "Both algorithms take O(m*n) time", 2)

# 6) 2*mn -> 2*(m*n)
$d.Content.Find.Execute("actually takes 2*mn time", $false, $false, $false, $false, $false, $true, 1, $false, "actually takes 2*(m*n) time", 2)

# 7) basic implementation takes mn time -> basic implementation takes m*n time
$d.Content.Find.Execute("basic implementation takes mn time", $false, $false, $false, $false, $false, $true, 1, $false, "basic implementation takes m*n time", 2)
